$d = $word.ActiveDocument

# Remove the two leading bullet-point paragraphs:
#   "Short description of what was done on project"
#   "Interesting items"
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()

# The "_GoBack" bookmark used to sit at the very end of the "This was
# chosen..." paragraph; move it to the very start of the document
# (now the "HVAC Helper" title paragraph), as a collapsed bookmark.
#
# Adding a bookmark exactly at document position 0 can cause its end
# marker to drift into the next paragraph, so we anchor it using a
# temporary marker character that is removed afterwards.
$anchor = $d.Range(0, 0)
$anchor.InsertBefore("x")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()
